$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column stays text (values like "22.163.86" or "0.9987" must
# not be auto-converted to numbers/dates by Excel's smart entry).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.163.86"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.560.41"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "0.9985"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "288.84"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.3800"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("D8").Value = "0.3302"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "43.92"
$ws.Range("E9").Value = "  -9.04%  "
$ws.Range("D10").Value = "1.149"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "0.07393"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "20.23"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "5.846"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "6.872"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "1.567.72"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "86.22"
$ws.Range("D19").Value = "0.06631"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").Value = "6.420"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "16.15"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "22.161.62"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "2.307"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "2.536"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").Value = "150.11"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "19.19"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").Value = "4.931"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "122.01"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "1.742.15"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "1.081"
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("D33").Value = "5.963"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").Value = "1.850"
$ws.Range("E34").Value = "  -8.03%  "
$ws.Range("D35").Value = "0.08264"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").Value = "9.353"
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("D37").Value = "0.02352"
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("D38").Value = "5.328"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "0.06267"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "0.2174"
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("D41").Value = "1.255"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "11.11"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "0.6091"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").Value = "0.9984"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "13.83"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "3.746"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "0.5906"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "1.998"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "122.41"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("D50").Value = "1.179"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("D51").Value = "0.07035"
$ws.Range("E51").Value = "  -2.75%  "
